$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - two newly-handed-off source files
# (21fbb4bd-f157-46e1-b5c7-6774a43be4b2 and b6dece24-d5b1-4107-b058-37915c45e933)
# are appended as row 4 / row 5 on each of the three sheets.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$file1 = "21fbb4bd-f157-46e1-b5c7-6774a43be4b2"
$file2 = "b6dece24-d5b1-4107-b058-37915c45e933"

$hash1 = "5542f86d38b1bd8f63e1dbbdeb8293c4ae18e192"
$hash2 = "2757c0915ebc045ce2ede704a80792c721e7ae07"

$md1 = $file1 + ".md"
$md2 = $file2 + ".md"

$xlfZh1 = $file1 + "." + $hash1 + ".zh-cn.xlf"
$xlfZh2 = $file2 + "." + $hash2 + ".zh-cn.xlf"
$xlfDe1 = $file1 + "." + $hash1 + ".de-de.xlf"
$xlfDe2 = $file2 + "." + $hash2 + ".de-de.xlf"

$statusReady = "Ready for handoff"
$extMd = ".md"
$include = "Include"
$epoch = "0001-01-01 00:00:00"

$dtOverview = "2016-03-23 06:03:20"
$dtZh = "2016-03-23 06:03:12"
$dtDe = "2016-03-23 06:03:20"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/9ff938606e9db46494ec6e23161c7418a92feb8a/e2e/"
$zhXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/29c92f6904eb1cb19a9934937642629a45e4a190/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
$deXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/428195fa63241df58dede60aff17a332d7919333/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (md link), B (zh-cn status), C (de-de status), D (date)
# ---------------------------------------------------------------------------

$ws1.Hyperlinks.Add($ws1.Range("A4"), ($mdUrlBase + $md1), "", "", $md1)
$ws1.Range("B4").Value = $statusReady
$ws1.Range("C4").Value = $statusReady
$ws1.Range("D4").Value = $dtOverview

$ws1.Hyperlinks.Add($ws1.Range("A5"), ($mdUrlBase + $md2), "", "", $md2)
$ws1.Range("B5").Value = $statusReady
$ws1.Range("C5").Value = $statusReady
$ws1.Range("D5").Value = $dtOverview

# ---------------------------------------------------------------------------
# Sheet "zh-cn": A (md link), B (ext), C (status), D (xlf link), E (handoff dt),
#                H (handback dt), J (reason)
# ---------------------------------------------------------------------------

$ws2.Hyperlinks.Add($ws2.Range("A4"), ($mdUrlBase + $md1), "", "", $md1)
$ws2.Range("B4").Value = $extMd
$ws2.Range("C4").Value = $statusReady
$ws2.Hyperlinks.Add($ws2.Range("D4"), ($zhXlfUrlBase + $xlfZh1), "", "", $xlfZh1)
$ws2.Range("E4").Value = $dtZh
$ws2.Range("H4").Value = $epoch
$ws2.Range("J4").Value = $include

$ws2.Hyperlinks.Add($ws2.Range("A5"), ($mdUrlBase + $md2), "", "", $md2)
$ws2.Range("B5").Value = $extMd
$ws2.Range("C5").Value = $statusReady
$ws2.Hyperlinks.Add($ws2.Range("D5"), ($zhXlfUrlBase + $xlfZh2), "", "", $xlfZh2)
$ws2.Range("E5").Value = $dtZh
$ws2.Range("H5").Value = $epoch
$ws2.Range("J5").Value = $include

# ---------------------------------------------------------------------------
# Sheet "de-de": A (md link), B (ext), C (status), D (xlf link), E (handoff dt),
#                H (handback dt), J (reason)
# ---------------------------------------------------------------------------

$ws3.Hyperlinks.Add($ws3.Range("A4"), ($mdUrlBase + $md1), "", "", $md1)
$ws3.Range("B4").Value = $extMd
$ws3.Range("C4").Value = $statusReady
$ws3.Hyperlinks.Add($ws3.Range("D4"), ($deXlfUrlBase + $xlfDe1), "", "", $xlfDe1)
$ws3.Range("E4").Value = $dtDe
$ws3.Range("H4").Value = $epoch
$ws3.Range("J4").Value = $include

$ws3.Hyperlinks.Add($ws3.Range("A5"), ($mdUrlBase + $md2), "", "", $md2)
$ws3.Range("B5").Value = $extMd
$ws3.Range("C5").Value = $statusReady
$ws3.Hyperlinks.Add($ws3.Range("D5"), ($deXlfUrlBase + $xlfDe2), "", "", $xlfDe2)
$ws3.Range("E5").Value = $dtDe
$ws3.Range("H5").Value = $epoch
$ws3.Range("J5").Value = $include

Write-Host "Report generated for handoff: rows 4-5 added to Overview, zh-cn, de-de"
